$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 56

# Columns A-D hold text-like values ("2023-06-18", "21:46:21", "Sunday", "25").
# Force them to be stored as text (matching the rest of the column) instead of
# letting Excel auto-convert them to a date serial / time / number.
$ws.Range("A$row").NumberFormat = "@"
$ws.Range("B$row").NumberFormat = "@"
$ws.Range("C$row").NumberFormat = "@"
$ws.Range("D$row").NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-18"
$ws.Cells.Item($row, 2).Value = "21:46:21"
$ws.Cells.Item($row, 3).Value = "Sunday"
$ws.Cells.Item($row, 4).Value = "25"

# Drop the temporary Text number-format again so the new row doesn't end up
# with explicit cell styling that the rest of the sheet doesn't have.
$ws.Range("A$row").ClearFormats()
$ws.Range("B$row").ClearFormats()
$ws.Range("C$row").ClearFormats()
$ws.Range("D$row").ClearFormats()

$ws.Cells.Item($row, 5).Value = 122092
$ws.Cells.Item($row, 6).Value = 133580
$ws.Cells.Item($row, 7).Value = 162532
$ws.Cells.Item($row, 8).Value = 133045
$ws.Cells.Item($row, 9).Value = 177394
$ws.Cells.Item($row, 10).Value = 115089
$ws.Cells.Item($row, 11).Value = 201552
$ws.Cells.Item($row, 12).Value = 225323
$ws.Cells.Item($row, 13).Value = 175460
$ws.Cells.Item($row, 14).Value = 103679
$ws.Cells.Item($row, 15).Value = 39253
$ws.Cells.Item($row, 16).Value = 33985
$ws.Cells.Item($row, 17).Value = 51849
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 36815
$ws.Cells.Item($row, 20).Value = -1
